$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.535.19'
$ws.Range('E2').Value = '  +2.42%  '

$ws.Range('D3').Value = '2.412.10'
$ws.Range('E3').Value = '  +8.53%  '

$ws.Range('E4').Value = '  -0.17%  '

$ws.Range('D5').Value = '''322.84'
$ws.Range('E5').Value = '  +11.84%  '

$ws.Range('D6').Value = '''103.83'
$ws.Range('E6').Value = '  -5.73%  '

$ws.Range('E7').Value = '  +3.55%  '

$ws.Range('E8').Value = '  -0.19%  '

$ws.Range('D9').Value = '''0.652'
$ws.Range('E9').Value = '  +9.02%  '

$ws.Range('D10').Value = '''41.93'
$ws.Range('E10').Value = '  -3.35%  '

$ws.Range('D11').Value = '''0.0946'
$ws.Range('E11').Value = '  +4.02%  '

$ws.Range('D12').Value = '''8.68'
$ws.Range('E12').Value = '  +0.66%  '

$ws.Range('E13').Value = '  +1.75%  '

$ws.Range('D14').Value = '''17.38'
$ws.Range('E14').Value = '  +17.05%  '

$ws.Range('E15').Value = '  +2.63%  '

$ws.Range('D16').Value = '2.776.78'

$ws.Range('D17').Value = '2.405.90'
$ws.Range('E17').Value = '  +7.84%  '

$ws.Range('D18').Value = '43.567.80'
$ws.Range('E18').Value = '  +2.86%  '

$ws.Range('E19').Value = '  +5.09%  '

$ws.Range('D20').Value = '''7.40'
$ws.Range('E20').Value = '  +3.31%  '

$ws.Range('E21').Value = '  +3.36%  '

$ws.Range('D22').Value = '''3.48'
$ws.Range('E22').Value = '  +3.79%  '

$ws.Range('D23').Value = '''260.91'
$ws.Range('E23').Value = '  +12.75%  '

$ws.Range('E24').Value = '  +1.72%  '

$ws.Range('D25').Value = '''9.68'
$ws.Range('E25').Value = '  +7.81%  '

$ws.Range('D26').Value = '''12.01'
$ws.Range('E26').Value = '  +5.12%  '

$ws.Range('D27').Value = '''0.999'
$ws.Range('E27').Value = '  -0.07%  '

$ws.Range('D28').Value = '''22.97'
$ws.Range('E28').Value = '  +10.09%  '

$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = '''2.24'
$ws.Range('E29').Value = '  +1.92%  '

$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').Value = '''179.16'
$ws.Range('E30').Value = '  +3.55%  '

$ws.Range('B31').Value = 'InjectiveProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D31').Value = '''38.21'
$ws.Range('E31').Value = '  +3.26%  '

$ws.Range('B32').Value = 'WEMIXToken'
$ws.Range('C32').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D32').Value = '''3.24'
$ws.Range('E32').Value = '  +1.43%  '

$ws.Range('D33').Value = '''0.0936'
$ws.Range('E33').Value = '  +7.22%  '

$ws.Range('D34').Value = '''5.98'
$ws.Range('E34').Value = '  +6.94%  '

$ws.Range('D35').Value = '''0.133'
$ws.Range('E35').Value = '  +5.31%  '

$ws.Range('D36').Value = '''4.91'
$ws.Range('E36').Value = '  -1.95%  '

$ws.Range('E37').Value = '  +0.58%  '

$ws.Range('E38').Value = '  -5.40%  '

$ws.Range('D39').Value = '''2.92'
$ws.Range('E39').Value = '  +22.01%  '

$ws.Range('E40').Value = '  +0.83%  '

$ws.Range('E41').Value = '  +25.11%  '

$ws.Range('D42').Value = '''0.234'
$ws.Range('E42').Value = '  +1.90%  '

$ws.Range('D43').Value = '''124.54'
$ws.Range('E43').Value = '  +22.55%  '

$ws.Range('D44').Value = '''69.11'
$ws.Range('E44').Value = '  -7.13%  '

$ws.Range('E45').Value = '  +0.27%  '

$ws.Range('D46').Value = '''12.66'

$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').Value = '''9.58'
$ws.Range('E47').Value = '  +13.27%  '

$ws.Range('B48').Value = 'THORChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D48').Value = '''5.68'
$ws.Range('E48').Value = '  +5.58%  '

$ws.Range('E49').Value = '  +3.13%  '

$ws.Range('D50').Value = '1.597.44'
$ws.Range('E50').Value = '  +13.18%  '

$ws.Range('E51').Value = '  +3.64%  '
